# Cambios simulador analista, Retanqueo Múltiple
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RetanqueoMultiple")

# Update data values (values are stored as literal quoted strings in this workbook)
$ws.Range("A2").Value = '"12962960"'
$ws.Range("Q2").Value = '"ROBERTO HERNAN"'

# Activate sheet, update selection and scroll the view so column K is leftmost
# (topLeftCell I1 -> K1, selection O9 -> P9)
$ws.Activate()
$ws.Range("P9").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 11
